# "Generate Report for Archive"
# The localization-status report is regenerated: the one file that was
# previously "Ready for handoff" has now moved on to "In Translation".
# That status string shows up in three places:
#   - Overview!E2  (zh-cn status column)
#   - Overview!F2  (de-de status column)
#   - zh-cn!C2     (Status column)
#   - de-de!C2     (Status column)
# After the text changes, the Status-ish columns are re-fit to their new,
# narrower content.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn (col E) and de-de (col F) status cells ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Columns("E:F").ColumnWidth = 12.5

# --- zh-cn sheet: Status column (col C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Columns("C:C").ColumnWidth = 12.5

# --- de-de sheet: Status column (col C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Columns("C:C").ColumnWidth = 12.5

Write-Output "Report regenerated: status updated to 'In Translation' and columns resized."
